# Update cryptocurrency price (D) and 1h volume change (E) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.948.73"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "1.647.64"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'309.81"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.3879"
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("D8").Value = "'0.3799"
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").Value = "'52.18"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'1.340"
$ws.Range("E10").Value = "  -4.60%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "'0.08437"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "'23.85"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").Value = "'7.042"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").Value = "'8.013"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "'0.00001306"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("D17").Value = "1.646.98"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'94.10"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "'0.06974"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'19.59"
$ws.Range("E20").Value = "  -5.15%  "
$ws.Range("D21").Value = "'6.934"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'13.73"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "23.952.86"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "'2.451"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").Value = "'2.942"
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").Value = "'153.38"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("D29").Value = "'5.411"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "'137.92"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("D31").Value = "'7.816"
$ws.Range("E31").Value = "  -4.20%  "
$ws.Range("D32").Value = "'2.510"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "1.830.51"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").Value = "  -5.79%  "
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").Value = "'6.742"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").Value = "'0.02926"
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("D38").Value = "'0.2671"
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("D39").Value = "'10.68"
$ws.Range("E39").Value = "  -4.53%  "
$ws.Range("D40").Value = "'0.09056"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").Value = "'13.36"
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("D43").Value = "'1.421"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'16.26"
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("D45").Value = "'0.6941"
$ws.Range("D46").Value = "'2.440"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("D47").Value = "'4.088"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "'0.08302"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "'134.26"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").Value = "'1.228"
$ws.Range("E51").Value = "  -4.42%  "
